# Generate Report for Handback
#
# Two files (6b5c63da-... and d677720c-...) have come back from
# localization "Handed back: in sync with en-US". This rewrites the
# Overview sheet and the two per-locale sheets (zh-cn, de-de) so that:
#   - the handed-back files are listed first, with their new status
#   - the per-locale sheets record the handback file/datetime in the
#     (until-now-empty) "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns (F, G, H)
#   - all hyperlinks point at the right targets for the new row order

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

# Drop the old hyperlinks before we shuffle rows around - Hyperlinks.Add
# never replaces an existing link on a cell, it just stacks another one
# on top, so any cell getting a fresh link needs its old one gone first.
$overview.Cells.Hyperlinks.Delete()

$overviewRows = @(
    @{ Uuid = "6b5c63da-fc13-41d9-bd68-060e19295db8"; Status = "Handed back: in sync with en-US"; Date = "2016-16-20 08:16:15" },
    @{ Uuid = "d677720c-cdb4-46d0-bc5e-3c22d6103b1d"; Status = "Handed back: in sync with en-US"; Date = "2016-16-20 08:16:15" },
    @{ Uuid = "4cafa085-c57d-468a-9fa4-8c2cc7ba367a"; Status = "In Translation";                  Date = "2016-14-20 08:14:55" },
    @{ Uuid = "95c62293-13e0-40e1-a42a-1384476e290e"; Status = "Ready for handoff";                Date = "2016-16-20 08:16:15" }
)

$mdUrl = @{
    "6b5c63da-fc13-41d9-bd68-060e19295db8" = "https://github.com/OpenLocalizationTest/oltest/blob/526d3d760458683629ce2225d0ccbf2c97dcc08b/e2e/6b5c63da-fc13-41d9-bd68-060e19295db8.md"
    "d677720c-cdb4-46d0-bc5e-3c22d6103b1d" = "https://github.com/OpenLocalizationTest/oltest/blob/526d3d760458683629ce2225d0ccbf2c97dcc08b/e2e/d677720c-cdb4-46d0-bc5e-3c22d6103b1d.md"
    "4cafa085-c57d-468a-9fa4-8c2cc7ba367a" = "https://github.com/OpenLocalizationTest/oltest/blob/07b16fb85c2c5fbde6acd0138a9ba63fdb08e0b5/e2e/4cafa085-c57d-468a-9fa4-8c2cc7ba367a.md"
    "95c62293-13e0-40e1-a42a-1384476e290e" = "https://github.com/OpenLocalizationTest/oltest/blob/141532668ce15081863b2d6f31ae5e05cc48f010/e2e/95c62293-13e0-40e1-a42a-1384476e290e.md"
}

$r = 2
foreach ($row in $overviewRows) {
    $fileName = "$($row.Uuid).md"

    $overview.Cells.Item($r, 1).Value = $fileName
    $overview.Cells.Item($r, 2).Value = $row.Status
    $overview.Cells.Item($r, 3).Value = $row.Status
    $overview.Cells.Item($r, 4).Value = $row.Date

    $overview.Hyperlinks.Add($overview.Cells.Item($r, 1), $mdUrl[$row.Uuid], "", "", $fileName)

    $r = $r + 1
}

# ---------------------------------------------------------------------
# Shared lookup tables used by both locale sheets
# ---------------------------------------------------------------------

$mdFileUrl = $mdUrl

$xlfUrl = @{
    "zh-cn" = @{
        "6b5c63da-fc13-41d9-bd68-060e19295db8" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5801ffcb7724553019b141c1db1d0585b989a6f5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/6b5c63da-fc13-41d9-bd68-060e19295db8.559078efbd17910c9a2f47a4733741b77076e371.zh-cn.xlf"
        "d677720c-cdb4-46d0-bc5e-3c22d6103b1d" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5801ffcb7724553019b141c1db1d0585b989a6f5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/d677720c-cdb4-46d0-bc5e-3c22d6103b1d.bf97de1068b6c9bda37af7d75c8e07e786819850.zh-cn.xlf"
        "4cafa085-c57d-468a-9fa4-8c2cc7ba367a" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6b94cb7df4020f6aa1f45b6174f6389a0fb7edbd/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4cafa085-c57d-468a-9fa4-8c2cc7ba367a.dfc21d39b86379109e2de3a9b3a82e481706ee2a.zh-cn.xlf"
        "95c62293-13e0-40e1-a42a-1384476e290e" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5801ffcb7724553019b141c1db1d0585b989a6f5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/95c62293-13e0-40e1-a42a-1384476e290e.858e23a74f0775ed4d0177f3b68a978874072582.zh-cn.xlf"
    }
    "de-de" = @{
        "6b5c63da-fc13-41d9-bd68-060e19295db8" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/68c0c6b098886cb125b0467f36a64b2e50997164/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/6b5c63da-fc13-41d9-bd68-060e19295db8.559078efbd17910c9a2f47a4733741b77076e371.de-de.xlf"
        "d677720c-cdb4-46d0-bc5e-3c22d6103b1d" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/68c0c6b098886cb125b0467f36a64b2e50997164/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/d677720c-cdb4-46d0-bc5e-3c22d6103b1d.bf97de1068b6c9bda37af7d75c8e07e786819850.de-de.xlf"
        "4cafa085-c57d-468a-9fa4-8c2cc7ba367a" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/38e708781318158b5fc05f1696d9d5c68c1cfb52/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4cafa085-c57d-468a-9fa4-8c2cc7ba367a.dfc21d39b86379109e2de3a9b3a82e481706ee2a.de-de.xlf"
        "95c62293-13e0-40e1-a42a-1384476e290e" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/68c0c6b098886cb125b0467f36a64b2e50997164/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/95c62293-13e0-40e1-a42a-1384476e290e.858e23a74f0775ed4d0177f3b68a978874072582.de-de.xlf"
    }
}

$xlfDisplay = @{
    "zh-cn" = @{
        "6b5c63da-fc13-41d9-bd68-060e19295db8" = "6b5c63da-fc13-41d9-bd68-060e19295db8.559078efbd17910c9a2f47a4733741b77076e371.zh-cn.xlf"
        "d677720c-cdb4-46d0-bc5e-3c22d6103b1d" = "d677720c-cdb4-46d0-bc5e-3c22d6103b1d.bf97de1068b6c9bda37af7d75c8e07e786819850.zh-cn.xlf"
        "4cafa085-c57d-468a-9fa4-8c2cc7ba367a" = "4cafa085-c57d-468a-9fa4-8c2cc7ba367a.dfc21d39b86379109e2de3a9b3a82e481706ee2a.zh-cn.xlf"
        "95c62293-13e0-40e1-a42a-1384476e290e" = "95c62293-13e0-40e1-a42a-1384476e290e.858e23a74f0775ed4d0177f3b68a978874072582.zh-cn.xlf"
    }
    "de-de" = @{
        "6b5c63da-fc13-41d9-bd68-060e19295db8" = "6b5c63da-fc13-41d9-bd68-060e19295db8.559078efbd17910c9a2f47a4733741b77076e371.de-de.xlf"
        "d677720c-cdb4-46d0-bc5e-3c22d6103b1d" = "d677720c-cdb4-46d0-bc5e-3c22d6103b1d.bf97de1068b6c9bda37af7d75c8e07e786819850.de-de.xlf"
        "4cafa085-c57d-468a-9fa4-8c2cc7ba367a" = "4cafa085-c57d-468a-9fa4-8c2cc7ba367a.dfc21d39b86379109e2de3a9b3a82e481706ee2a.de-de.xlf"
        "95c62293-13e0-40e1-a42a-1384476e290e" = "95c62293-13e0-40e1-a42a-1384476e290e.858e23a74f0775ed4d0177f3b68a978874072582.de-de.xlf"
    }
}

# Per-row data (shared shape across locale sheets; only dates/target
# file differ by locale, handled through the hashtables above).
$localeRows = @(
    @{ Uuid = "6b5c63da-fc13-41d9-bd68-060e19295db8"; Status = "Handed back: in sync with en-US"; HasHandback = $true },
    @{ Uuid = "d677720c-cdb4-46d0-bc5e-3c22d6103b1d"; Status = "Handed back: in sync with en-US"; HasHandback = $true },
    @{ Uuid = "4cafa085-c57d-468a-9fa4-8c2cc7ba367a"; Status = "In Translation";                  HasHandback = $false },
    @{ Uuid = "95c62293-13e0-40e1-a42a-1384476e290e"; Status = "Ready for handoff";                HasHandback = $false }
)

$handoffDateTime = @{
    "zh-cn" = @{
        "6b5c63da-fc13-41d9-bd68-060e19295db8" = "2016-03-20 08:16:11"
        "d677720c-cdb4-46d0-bc5e-3c22d6103b1d" = "2016-03-20 08:16:11"
        "4cafa085-c57d-468a-9fa4-8c2cc7ba367a" = "2016-03-20 08:14:52"
        "95c62293-13e0-40e1-a42a-1384476e290e" = "2016-03-20 08:16:11"
    }
    "de-de" = @{
        "6b5c63da-fc13-41d9-bd68-060e19295db8" = "2016-03-20 08:16:15"
        "d677720c-cdb4-46d0-bc5e-3c22d6103b1d" = "2016-03-20 08:16:15"
        "4cafa085-c57d-468a-9fa4-8c2cc7ba367a" = "2016-03-20 08:14:55"
        "95c62293-13e0-40e1-a42a-1384476e290e" = "2016-03-20 08:16:15"
    }
}

$handbackDateTime = @{
    "zh-cn" = @{
        "6b5c63da-fc13-41d9-bd68-060e19295db8" = "2016-03-20 08:16:29"
        "d677720c-cdb4-46d0-bc5e-3c22d6103b1d" = "2016-03-20 08:16:29"
        "4cafa085-c57d-468a-9fa4-8c2cc7ba367a" = "0001-01-01 00:00:00"
        "95c62293-13e0-40e1-a42a-1384476e290e" = "0001-01-01 00:00:00"
    }
    "de-de" = @{
        "6b5c63da-fc13-41d9-bd68-060e19295db8" = "2016-03-20 08:16:35"
        "d677720c-cdb4-46d0-bc5e-3c22d6103b1d" = "2016-03-20 08:16:35"
        "4cafa085-c57d-468a-9fa4-8c2cc7ba367a" = "0001-01-01 00:00:00"
        "95c62293-13e0-40e1-a42a-1384476e290e" = "0001-01-01 00:00:00"
    }
}

function Update-LocaleSheet($sheetName) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Cells.Hyperlinks.Delete()

    $r = 2
    foreach ($row in $localeRows) {
        $uuid = $row.Uuid
        $mdFileName = "$uuid.md"

        $ws.Cells.Item($r, 1).Value = $mdFileName          # A: Source File Name
        $ws.Cells.Item($r, 2).Value = ".md"                 # B: File Extension
        $ws.Cells.Item($r, 3).Value = $row.Status           # C: Status
        $ws.Cells.Item($r, 4).Value = $xlfDisplay[$sheetName][$uuid]   # D: Latest Handoff File
        $ws.Cells.Item($r, 5).Value = $handoffDateTime[$sheetName][$uuid]  # E: Latest Handoff Datetime

        $ws.Hyperlinks.Add($ws.Cells.Item($r, 1), $mdFileUrl[$uuid], "", "", $mdFileName)
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), $mdFileUrl[$uuid], "", "", ".md")
        $ws.Hyperlinks.Add($ws.Cells.Item($r, 4), $xlfUrl[$sheetName][$uuid], "", "", $xlfDisplay[$sheetName][$uuid])

        if ($row.HasHandback) {
            $ws.Cells.Item($r, 6).Value = $mdFileName                        # F: Latest Target File
            $ws.Cells.Item($r, 7).Value = $xlfDisplay[$sheetName][$uuid]     # G: Latest Handback File
            $ws.Cells.Item($r, 8).Value = $handbackDateTime[$sheetName][$uuid]  # H: Latest Handback DateTime

            $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $mdFileUrl[$uuid], "", "", $mdFileName)
            $ws.Hyperlinks.Add($ws.Cells.Item($r, 7), $xlfUrl[$sheetName][$uuid], "", "", $xlfDisplay[$sheetName][$uuid])
        } else {
            $ws.Cells.Item($r, 8).Value = $handbackDateTime[$sheetName][$uuid]  # H stays the "never" sentinel
        }

        $ws.Cells.Item($r, 9).Value = "Include"             # I: Handoff Reason

        $r = $r + 1
    }
}

Update-LocaleSheet "zh-cn"
Update-LocaleSheet "de-de"
